$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17106294631958
$ws.Range("B1").Value = 2.438135147094727
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.367327928543091
$ws.Range("E1").Value = 1.234786748886108
